# "limit volunteers, lidar com turnos extra"
#
# The volunteer roster had grown to Vol48 with a 40-slot availability
# array; this trims it back down to Vol30 (removing the Vol31..Vol48
# rows) and shrinks the availability array literal stored in column C
# from 40 entries to 28 entries for every remaining data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last 18 volunteer rows (Vol31..Vol48), rows 73-90.
$ws.Range("A73:C90").EntireRow.Delete()

# Shrink the "Disponibilidade" availability array from 40 to 28 slots
# for every remaining person (rows 2-72).
$newAvailability = "[1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1]"
$ws.Range("C2:C72").Value = $newAvailability

# Keep the header selection in sync with the shrunk data range and
# scroll the view back up towards the remaining rows.
$ws.Range("C2:C72").Select()
$excel.ActiveWindow.ScrollRow = 43
